$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, shifting existing rows 38:163 down to 39:164.
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44701
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100114007
$ws.Range("G38").Value = "Jengibre"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 40
$ws.Range("K38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = 20000
$ws.Range("N38").Value = "$/caja 13 kilos"
$ws.Range("O38").Value = "Perú"
$ws.Range("P38").Value = 1538
$ws.Range("Q38").Value = 13
$ws.Range("R38").Value = "Hortaliza"
